# ============================================================
# Update 合肥-漫展信息.xlsx: apply new/changed exhibition rows
# to both the '展览' (sheet 1) and '全部类型' (sheet 4) tabs,
# which carry identical data in this workbook.
# ============================================================

$wb = $excel.ActiveWorkbook

foreach ($sheetIndex in 1,4) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    # --- Simple numeric 'want to go' count bumps on existing rows ---
    $ws.Cells.Item(3,6).Value = 7606
    $ws.Cells.Item(6,6).Value = 32
    $ws.Cells.Item(9,6).Value = 5752
    $ws.Cells.Item(10,6).Value = 148
    $ws.Cells.Item(11,6).Value = 11
    $ws.Cells.Item(12,6).Value = 19
    $ws.Cells.Item(13,6).Value = 1760

    # --- Insert two new rows before row 17 so the sheet grows from
    #     17 used rows (1 header + 16 events) to 19 (1 header + 18 events).
    #     This pushes the former row 17 ('第二届漫画城市动漫展') down to row 19.
    $ws.Rows.Item(17).Insert()
    $ws.Rows.Item(17).Insert()

    # --- Column A on every data row carries a bold/bordered/centered
    #     style (the row-number style). Copy that formatting from row 16
    #     onto the two freshly inserted rows before filling them in.
    $ws.Cells.Item(16,1).Copy()
    $ws.Cells.Item(17,1).PasteSpecial(-4122)
    $ws.Cells.Item(16,1).Copy()
    $ws.Cells.Item(18,1).PasteSpecial(-4122)
    $ws.Application.CutCopyMode = $false

    # --- Rewrite rows 14-19 in full: the event list was reordered
    #     chronologically and two brand-new events were added, so the
    #     tail of the table (everything from 2024-02-17 onward) is
    #     re-entered with its final, post-edit content. ---

    # Row 14: 合肥·2024运动新春动漫庆典（全ip）
    $ws.Cells.Item(14,1).Value = 13
    $ws.Cells.Item(14,2).NumberFormat = "@"
    $ws.Cells.Item(14,2).Value = '2024-02-17'
    $ws.Cells.Item(14,2).Style = "Normal"
    $ws.Cells.Item(14,3).Value = '合肥·2024运动新春动漫庆典（全ip）'
    $ws.Cells.Item(14,4).Value = '锦绣大道与清潭路交口东北角 李宁体育公园'
    $ws.Cells.Item(14,5).Value = '2024.02.17 09:00-02.17 17:00'
    $ws.Cells.Item(14,6).Value = 1256
    $ws.Cells.Item(14,7).NumberFormat = "@"
    $ws.Cells.Item(14,7).Value = '65'
    $ws.Cells.Item(14,7).Style = "Normal"
    $ws.Cells.Item(14,8).Value = $false
    $ws.Cells.Item(14,9).Value = 'https://show.bilibili.com/platform/detail.html?id=79918'
    $ws.Cells.Item(14,10).Value = '//i0.hdslb.com/bfs/openplatform/202312/vzuMc0sJ1702902061660.jpeg'

    # Row 15: 合肥·安徽马娘only
    $ws.Cells.Item(15,1).Value = 14
    $ws.Cells.Item(15,2).NumberFormat = "@"
    $ws.Cells.Item(15,2).Value = '2024-02-19'
    $ws.Cells.Item(15,2).Style = "Normal"
    $ws.Cells.Item(15,3).Value = '合肥·安徽马娘only'
    $ws.Cells.Item(15,4).Value = '桐城路与庐江路交叉口西南80米 赤阑桥文玩大厦'
    $ws.Cells.Item(15,5).Value = '2024.02.19 09:00-02.19 17:00'
    $ws.Cells.Item(15,6).Value = 279
    $ws.Cells.Item(15,7).NumberFormat = "@"
    $ws.Cells.Item(15,7).Value = '68'
    $ws.Cells.Item(15,7).Style = "Normal"
    $ws.Cells.Item(15,8).Value = $false
    $ws.Cells.Item(15,9).Value = 'https://show.bilibili.com/platform/detail.html?id=78286'
    $ws.Cells.Item(15,10).Value = '//i1.hdslb.com/bfs/openplatform/202311/721L5pIZ1699428443216.jpeg'

    # Row 16: 合肥·星芒1.5动漫嘉年华
    $ws.Cells.Item(16,1).Value = 15
    $ws.Cells.Item(16,2).NumberFormat = "@"
    $ws.Cells.Item(16,2).Value = '2024-03-02'
    $ws.Cells.Item(16,2).Style = "Normal"
    $ws.Cells.Item(16,3).Value = '合肥·星芒1.5动漫嘉年华'
    $ws.Cells.Item(16,4).Value = '山西路与太原路交叉口 挥动体育'
    $ws.Cells.Item(16,5).Value = '2024.03.02 09:30-03.02 17:30'
    $ws.Cells.Item(16,6).Value = 24
    $ws.Cells.Item(16,7).NumberFormat = "@"
    $ws.Cells.Item(16,7).Value = '55'
    $ws.Cells.Item(16,7).Style = "Normal"
    $ws.Cells.Item(16,8).Value = $false
    $ws.Cells.Item(16,9).Value = 'https://show.bilibili.com/platform/detail.html?id=81267'
    $ws.Cells.Item(16,10).Value = '//i0.hdslb.com/bfs/openplatform/202401/GWidiefU1706003134747.jpeg'

    # Row 17: 合肥·CW国潮动漫游戏嘉年华
    $ws.Cells.Item(17,1).Value = 16
    $ws.Cells.Item(17,2).NumberFormat = "@"
    $ws.Cells.Item(17,2).Value = '2024-03-16'
    $ws.Cells.Item(17,2).Style = "Normal"
    $ws.Cells.Item(17,3).Value = '合肥·CW国潮动漫游戏嘉年华'
    $ws.Cells.Item(17,4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
    $ws.Cells.Item(17,5).Value = '2024.03.16 09:30-03.17 17:00'
    $ws.Cells.Item(17,6).Value = 2
    $ws.Cells.Item(17,7).Value = '不可售'
    $ws.Cells.Item(17,8).Value = $true
    $ws.Cells.Item(17,9).Value = 'https://show.bilibili.com/platform/detail.html?id=81284'
    $ws.Cells.Item(17,10).Value = '//i0.hdslb.com/bfs/openplatform/202401/38B92fWF1705995243803.jpeg'

    # Row 18: 合肥· 第二届漫画城市动漫展 -故事再次开始
    $ws.Cells.Item(18,1).Value = 17
    $ws.Cells.Item(18,2).NumberFormat = "@"
    $ws.Cells.Item(18,2).Value = '2024-04-04'
    $ws.Cells.Item(18,2).Style = "Normal"
    $ws.Cells.Item(18,3).Value = '合肥· 第二届漫画城市动漫展 -故事再次开始'
    $ws.Cells.Item(18,4).Value = '凤淮路与固镇路西北角 庐阳全民健身中心'
    $ws.Cells.Item(18,5).Value = '2024.04.04 09:00-04.05 17:00'
    $ws.Cells.Item(18,6).Value = 5506
    $ws.Cells.Item(18,7).NumberFormat = "@"
    $ws.Cells.Item(18,7).Value = '60'
    $ws.Cells.Item(18,7).Style = "Normal"
    $ws.Cells.Item(18,8).Value = $false
    $ws.Cells.Item(18,9).Value = 'https://show.bilibili.com/platform/detail.html?id=78898'
    $ws.Cells.Item(18,10).Value = '//i2.hdslb.com/bfs/openplatform/202311/244eBWip1700711342120.jpeg'

    # Row 19: 合肥·梦时空SPO1动漫展
    $ws.Cells.Item(19,1).Value = 18
    $ws.Cells.Item(19,2).NumberFormat = "@"
    $ws.Cells.Item(19,2).Value = '2024-05-18'
    $ws.Cells.Item(19,2).Style = "Normal"
    $ws.Cells.Item(19,3).Value = '合肥·梦时空SPO1动漫展'
    $ws.Cells.Item(19,4).Value = '阜阳路16号 银瑞林国际大酒店'
    $ws.Cells.Item(19,5).Value = '2024.05.18 10:00-05.18 17:00'
    $ws.Cells.Item(19,6).Value = 61
    $ws.Cells.Item(19,7).NumberFormat = "@"
    $ws.Cells.Item(19,7).Value = '60'
    $ws.Cells.Item(19,7).Style = "Normal"
    $ws.Cells.Item(19,8).Value = $false
    $ws.Cells.Item(19,9).Value = 'https://show.bilibili.com/platform/detail.html?id=80207'
    $ws.Cells.Item(19,10).Value = '//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg'
}

Write-Output "done"